# Apply the "add three more unsupervised methods" edit to the
# "Other metrics(unlabeled pctg, cluster num...)" sheet (4th sheet).
#
# Net effect (per the target diff):
#   - The "cluster_num" metric row is dropped for every dataset, leaving
#     only "unlabeled_pctg" and "pred_type_max_pctg" rows per dataset
#     (3 rows/dataset -> 2 rows/dataset, rows 4:33 -> rows 4:23).
#   - The now-unused "cluster_num" shared string disappears.
#   - The sheet title cell (A1) text shortens from
#     "Other metrics(unlabeled pctg, cluster num...)" to
#     "Other metrics(unlabeled pctg...)".
#   - The sheet tab itself is renamed from
#     "Other metrics(unlabeled pctg, c" to "Supervised methods other metric".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# Each of the 10 datasets occupies a 3-row block starting at row 4:
#   row  4 -> unlabeled_pctg      (dataset 1)
#   row  5 -> cluster_num         (dataset 1)   <- delete
#   row  6 -> pred_type_max_pctg  (dataset 1)
#   row  7 -> unlabeled_pctg      (dataset 2)
#   ...
# The "cluster_num" row is the middle row of every block: 5, 8, 11, ...,
# 32. Delete from the bottom up so earlier row numbers stay valid as we go.
$rowsToDelete = @(32, 29, 26, 23, 20, 17, 14, 11, 8, 5)
foreach ($r in $rowsToDelete) {
    $ws.Rows($r).Delete()
}

# Shorten the sheet's title cell text (row 1).
$ws.Range("A1").Value = "Other metrics(unlabeled pctg...)"

# Rename the sheet tab.
$ws.Name = "Supervised methods other metric"
